# demo_data/timeoftips.xlsx -- "making a lot of changes, adding instructions"
#
# Renames the worksheet, relabels the rainfall-value column header, moves the
# active selection, widens column B to fit the new header text, and nudges the
# saved window position, matching the authored commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename "Sheet1" -> "rainfall_data"
$ws.Name = "rainfall_data"

# Relabel the header in B1 ("r3ba1r01r" -> "rainfall_value")
$ws.Range("B1").Value = "rainfall_value"

# Widen column B so the longer header text fits (was bestFit width 11, now ~13.29)
$ws.Columns.Item(2).ColumnWidth = 12.5

# Move the current selection from H15 to C21
$ws.Range("C21").Select() | Out-Null

# Reposition the saved window location
$aw = $excel.Windows.Item(1)
$aw.Left = -28920
$aw.Top = -120
